$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Chee Tian" evaluator scores for the s2s_lstm_att_long_win
# model (columns AA:AF) on rows 13-17 -- every score is 1.
$ws.Range("AA13:AF17").Value = 1

# Leave the selection where the user finished editing: cell AF18.
$ws.Range("AF18").Select()
